# Unit-4 Authorised Assignment Brief - fix swapped logo "name" identifiers
# on the header/footer inline pictures (Pearson logo + BTEC logo).
#
# The three affected pictures all live in headers/footers, not the main
# story, so we have to reach them through Sections(1).Headers/Footers
# rather than ActiveDocument.InlineShapes.
#
# InlineShape has no writable .Name property when addressed straight off
# a Range (that mirrors real Word: inline pictures don't expose Name in
# the object model the way floating Shapes do). Selecting the picture
# first and then going through Selection.InlineShapes(1).Name = "..." is
# the COM-automation pattern that actually lands the rename.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Primary (default) footer: Pearson Edexcel logo -----------------
# docPr/cNvPr id="2" -> rename image2.png to image1.png
$footerPrimary = $sec.Footers.Item(1)
$pearsonPrimary = $footerPrimary.Range.InlineShapes.Item(1)
[void]$pearsonPrimary.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

# --- First-page header: BTEC logo ------------------------------------
# docPr/cNvPr id="1" -> rename image1.jpg to image2.jpg
$headerFirstPage = $sec.Headers.Item(2)
$btecLogo = $headerFirstPage.Range.InlineShapes.Item(1)
[void]$btecLogo.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.jpg"

# --- First-page footer: Pearson Edexcel logo -------------------------
# docPr/cNvPr id="3" -> rename image2.png to image1.png
$footerFirstPage = $sec.Footers.Item(2)
$pearsonFirstPage = $footerFirstPage.Range.InlineShapes.Item(1)
[void]$pearsonFirstPage.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

Write-Output "Renamed header/footer logo inline shapes."
